$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new value for L5 (Diagnoses Expansion time for the 1,000,000 cohort row)
$ws.Range("L5").Value = "14479 seconds"

# Update the view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("L5").Select()
